$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Rtn4"
$ws.Cells.Item(2, 3).Value2 = "Rtn4r"
$ws.Cells.Item(2, 4).Value2 = "FAPs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 43.96295866666666
$ws.Cells.Item(2, 8).Value2 = 131.888876
$ws.Cells.Item(2, 9).Value2 = 0.1193823486802574
$ws.Cells.Item(2, 10).Value2 = 0.1297146081693155
$ws.Cells.Item(2, 11).Value2 = 1
$ws.Cells.Item(2, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 13).Value2 = 0.2214103333333333
$ws.Cells.Item(2, 14).Value2 = 0.664231
$ws.Cells.Item(2, 15).Value2 = 0.5896903116237344
$ws.Cells.Item(2, 16).Value2 = 0.6831206522767569
$ws.Cells.Item(2, 17).Value2 = 9.733853332706222
$ws.Cells.Item(2, 18).Value2 = 87.60467999435599
$ws.Cells.Item(2, 19).Value2 = 0.07039861439563429
$ws.Cells.Item(2, 20).Value2 = 0.08861072774244672

$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Rtn4"
$ws.Cells.Item(3, 3).Value2 = "Rtn4r"
$ws.Cells.Item(3, 4).Value2 = "MuSCs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 43.96295866666666
$ws.Cells.Item(3, 8).Value2 = 131.888876
$ws.Cells.Item(3, 9).Value2 = 0.1193823486802574
$ws.Cells.Item(3, 10).Value2 = 0.1297146081693155
$ws.Cells.Item(3, 11).Value2 = 1
$ws.Cells.Item(3, 12).Value2 = 0.5
$ws.Cells.Item(3, 13).Value2 = 0.1540585
$ws.Cells.Item(3, 14).Value2 = 0.308117
$ws.Cells.Item(3, 15).Value2 = 0.4103096883762655
$ws.Cells.Item(3, 16).Value2 = 0.3168793477232431
$ws.Cells.Item(3, 17).Value2 = 6.772867467748664
$ws.Cells.Item(3, 18).Value2 = 40.63720480649199
$ws.Cells.Item(3, 19).Value2 = 0.04898373428462308
$ws.Cells.Item(3, 20).Value2 = 0.04110388042686874

$ws.Cells.Item(4, 1).Value2 = "FAPs"
$ws.Cells.Item(4, 2).Value2 = "Rtn4"
$ws.Cells.Item(4, 3).Value2 = "Rtn4r"
$ws.Cells.Item(4, 4).Value2 = "FAPs"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 72.02213166666667
$ws.Cells.Item(4, 8).Value2 = 216.066395
$ws.Cells.Item(4, 9).Value2 = 0.1955776293519722
$ws.Cells.Item(4, 10).Value2 = 0.212504409894141
$ws.Cells.Item(4, 11).Value2 = 1
$ws.Cells.Item(4, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(4, 13).Value2 = 0.2214103333333333
$ws.Cells.Item(4, 14).Value2 = 0.664231
$ws.Cells.Item(4, 15).Value2 = 0.5896903116237344
$ws.Cells.Item(4, 16).Value2 = 0.6831206522767569
$ws.Cells.Item(4, 17).Value2 = 15.94644417969389
$ws.Cells.Item(4, 18).Value2 = 143.517997617245
$ws.Cells.Item(4, 19).Value2 = 0.1153302331991957
$ws.Cells.Item(4, 20).Value2 = 0.1451661510985729

$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Rtn4"
$ws.Cells.Item(5, 3).Value2 = "Rtn4r"
$ws.Cells.Item(5, 4).Value2 = "MuSCs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 72.02213166666667
$ws.Cells.Item(5, 8).Value2 = 216.066395
$ws.Cells.Item(5, 9).Value2 = 0.1955776293519722
$ws.Cells.Item(5, 10).Value2 = 0.212504409894141
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0.5
$ws.Cells.Item(5, 13).Value2 = 0.1540585
$ws.Cells.Item(5, 14).Value2 = 0.308117
$ws.Cells.Item(5, 15).Value2 = 0.4103096883762655
$ws.Cells.Item(5, 16).Value2 = 0.3168793477232431
$ws.Cells.Item(5, 17).Value2 = 11.09562157136917
$ws.Cells.Item(5, 18).Value2 = 66.57372942821499
$ws.Cells.Item(5, 19).Value2 = 0.08024739615277647
$ws.Cells.Item(5, 20).Value2 = 0.06733825879556811

$ws.Cells.Item(6, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(6, 2).Value2 = "Rtn4"
$ws.Cells.Item(6, 3).Value2 = "Rtn4r"
$ws.Cells.Item(6, 4).Value2 = "FAPs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 75.27587666666666
$ws.Cells.Item(6, 8).Value2 = 225.82763
$ws.Cells.Item(6, 9).Value2 = 0.2044132430569516
$ws.Cells.Item(6, 10).Value2 = 0.2221047250357578
$ws.Cells.Item(6, 11).Value2 = 1
$ws.Cells.Item(6, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(6, 13).Value2 = 0.2214103333333333
$ws.Cells.Item(6, 14).Value2 = 0.664231
$ws.Cells.Item(6, 15).Value2 = 0.5896903116237344
$ws.Cells.Item(6, 16).Value2 = 0.6831206522767569
$ws.Cells.Item(6, 17).Value2 = 16.66685694472556
$ws.Cells.Item(6, 18).Value2 = 150.00171250253
$ws.Cells.Item(6, 19).Value2 = 0.120540508998272
$ws.Cells.Item(6, 20).Value2 = 0.1517243246401766

$ws.Cells.Item(7, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(7, 2).Value2 = "Rtn4"
$ws.Cells.Item(7, 3).Value2 = "Rtn4r"
$ws.Cells.Item(7, 4).Value2 = "MuSCs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 75.27587666666666
$ws.Cells.Item(7, 8).Value2 = 225.82763
$ws.Cells.Item(7, 9).Value2 = 0.2044132430569516
$ws.Cells.Item(7, 10).Value2 = 0.2221047250357578
$ws.Cells.Item(7, 11).Value2 = 1
$ws.Cells.Item(7, 12).Value2 = 0.5
$ws.Cells.Item(7, 13).Value2 = 0.1540585
$ws.Cells.Item(7, 14).Value2 = 0.308117
$ws.Cells.Item(7, 15).Value2 = 0.4103096883762655
$ws.Cells.Item(7, 16).Value2 = 0.3168793477232431
$ws.Cells.Item(7, 17).Value2 = 11.59688864545167
$ws.Cells.Item(7, 18).Value2 = 69.58133187271
$ws.Cells.Item(7, 19).Value2 = 0.08387273405867965
$ws.Cells.Item(7, 20).Value2 = 0.07038040039558119

$ws.Cells.Item(8, 1).Value2 = "MuSCs"
$ws.Cells.Item(8, 2).Value2 = "Rtn4"
$ws.Cells.Item(8, 3).Value2 = "Rtn4r"
$ws.Cells.Item(8, 4).Value2 = "FAPs"
$ws.Cells.Item(8, 5).Value2 = 2
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 87.99833699999999
$ws.Cells.Item(8, 8).Value2 = 175.996674
$ws.Cells.Item(8, 9).Value2 = 0.2389613545046087
$ws.Cells.Item(8, 10).Value2 = 0.1730952624618072
$ws.Cells.Item(8, 11).Value2 = 1
$ws.Cells.Item(8, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(8, 13).Value2 = 0.2214103333333333
$ws.Cells.Item(8, 14).Value2 = 0.664231
$ws.Cells.Item(8, 15).Value2 = 0.5896903116237344
$ws.Cells.Item(8, 16).Value2 = 0.6831206522767569
$ws.Cells.Item(8, 17).Value2 = 19.483741127949
$ws.Cells.Item(8, 18).Value2 = 116.902446767694
$ws.Cells.Item(8, 19).Value2 = 0.1409131956038524
$ws.Cells.Item(8, 20).Value2 = 0.1182449485989262

$ws.Cells.Item(9, 1).Value2 = "MuSCs"
$ws.Cells.Item(9, 2).Value2 = "Rtn4"
$ws.Cells.Item(9, 3).Value2 = "Rtn4r"
$ws.Cells.Item(9, 4).Value2 = "MuSCs"
$ws.Cells.Item(9, 5).Value2 = 2
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 87.99833699999999
$ws.Cells.Item(9, 8).Value2 = 175.996674
$ws.Cells.Item(9, 9).Value2 = 0.2389613545046087
$ws.Cells.Item(9, 10).Value2 = 0.1730952624618072
$ws.Cells.Item(9, 11).Value2 = 1
$ws.Cells.Item(9, 12).Value2 = 0.5
$ws.Cells.Item(9, 13).Value2 = 0.1540585
$ws.Cells.Item(9, 14).Value2 = 0.308117
$ws.Cells.Item(9, 15).Value2 = 0.4103096883762655
$ws.Cells.Item(9, 16).Value2 = 0.3168793477232431
$ws.Cells.Item(9, 17).Value2 = 13.5568918007145
$ws.Cells.Item(9, 18).Value2 = 54.22756720285799
$ws.Cells.Item(9, 19).Value2 = 0.09804815890075633
$ws.Cells.Item(9, 20).Value2 = 0.05485031386288104

$ws.Cells.Item(10, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(10, 2).Value2 = "Rtn4"
$ws.Cells.Item(10, 3).Value2 = "Rtn4r"
$ws.Cells.Item(10, 4).Value2 = "FAPs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 88.99412
$ws.Cells.Item(10, 8).Value2 = 266.98236
$ws.Cells.Item(10, 9).Value2 = 0.24166542440621
$ws.Cells.Item(10, 10).Value2 = 0.2625809944389785
$ws.Cells.Item(10, 11).Value2 = 1
$ws.Cells.Item(10, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(10, 13).Value2 = 0.2214103333333333
$ws.Cells.Item(10, 14).Value2 = 0.664231
$ws.Cells.Item(10, 15).Value2 = 0.5896903116237344
$ws.Cells.Item(10, 16).Value2 = 0.6831206522767569
$ws.Cells.Item(10, 17).Value2 = 19.70421777390667
$ws.Cells.Item(10, 18).Value2 = 177.33795996516
$ws.Cells.Item(10, 19).Value2 = 0.14250775942678
$ws.Cells.Item(10, 20).Value2 = 0.1793745001966345

$ws.Cells.Item(11, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value2 = "Rtn4"
$ws.Cells.Item(11, 3).Value2 = "Rtn4r"
$ws.Cells.Item(11, 4).Value2 = "MuSCs"
$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 88.99412
$ws.Cells.Item(11, 8).Value2 = 266.98236
$ws.Cells.Item(11, 9).Value2 = 0.24166542440621
$ws.Cells.Item(11, 10).Value2 = 0.2625809944389785
$ws.Cells.Item(11, 11).Value2 = 1
$ws.Cells.Item(11, 12).Value2 = 0.5
$ws.Cells.Item(11, 13).Value2 = 0.1540585
$ws.Cells.Item(11, 14).Value2 = 0.308117
$ws.Cells.Item(11, 15).Value2 = 0.4103096883762655
$ws.Cells.Item(11, 16).Value2 = 0.3168793477232431
$ws.Cells.Item(11, 17).Value2 = 13.71030063602
$ws.Cells.Item(11, 18).Value2 = 82.26180381611998
$ws.Cells.Item(11, 19).Value2 = 0.09915766497942997
$ws.Cells.Item(11, 20).Value2 = 0.08320649424234403
